# fixing 2D x, y weirdness
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the wording of the existing comment (C3): "boundary condition" -> "boundary_conditions"
$ws.Range("C3").Value = "adding @njit to the boundary_conditions function"

# Fill in the new run data for rows 4-6 (A/B/C), and advance the selection to A8

# Row 4
$ws.Range("A4").Formula = "=(123.54788+125.00538+115.54499)/3"
$ws.Range("C4").Value = "removing @njit from the boundary_conditions function, removing cons2prim in w_half and get_flux, and adding @njit to w_half, get_flux, right_eigenvectors, and left_eigenvectors"

# Row 5
$ws.Range("A5").Formula = "=(12.31002+12.6224+13.26979)/3"
$ws.Range("C5").Value = "adding @njit to weno and lf_flux"

# Row 6
$ws.Range("A6").Value = 468.19565
$ws.Range("C6").Value = " 128x128 grid, CFL = 0.1, t1 = 0.3, case = 0"

$ws.Range("A8").Select()

$wb.Save()
